# Apply the "Add updated evaluation results + initial version sentence
# transformers learning" edit to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string / label edits ---------------------------------------
# "Score" sub-header (shared by all four metric groups) -> "score"
$ws.Range("B2").Value = "score"
$ws.Range("D2").Value = "score"
$ws.Range("F2").Value = "score"
$ws.Range("H2").Value = "score"

# Row labels for the two lower rows get annotated with their values.
$ws.Range("A5").Value = "HunFlair (1.030)"
$ws.Range("A6").Value = "UMLS (1.375)"

# --- Border/formatting normalisation ------------------------------------
# D2/H2 pick up the same left-hair border that B2/F2 already use.
$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# Row 6 (the "UMLS (1.375)" row) had D/E/H/I using a border style that
# didn't match rows 3 and 5 - line it up with those rows.
$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("I3").Copy()
$ws.Range("I6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Updated evaluation numbers for row 6 -------------------------------
$ws.Range("B6").Value = 51.39328
$ws.Range("C6").Value = 2.65853227770513
$ws.Range("D6").Value = 33.45576
$ws.Range("E6").Value = 2.79199697356569
$ws.Range("F6").Value = 48.2352
$ws.Range("G6").Value = 2.86973771379895
$ws.Range("H6").Value = 48.22435
$ws.Range("I6").Value = 2.82599274140964

# --- Selection moves to H6 ----------------------------------------------
[void]$ws.Range("H6").Select()
